$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows of data
$ws.Range("A12").Value = "mexerica ponkan"
$ws.Range("B12").Value = "x"
$ws.Range("C12").Value = "x"

$ws.Range("A13").Value = "hortelã"
$ws.Range("B13").Value = "x"

# Column A needs to widen to fit the new longer entry ("mexerica ponkan")
$ws.Columns("A:A").ColumnWidth = 15.5

# Update selection to match the target state
$ws.Range("F26").Select()
